$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Add a new row 16 whose format mirrors the existing data rows (row 15),
#    then shift the tail of the table (old rows 12-15) down by one row to
#    make room for the new "Jurisdiction" row at position 12.
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)

$ws.Range("A16").Value = $ws.Range("A15").Value()
$ws.Range("B16").Value = $ws.Range("B15").Value()
$ws.Range("A15").Value = $ws.Range("A14").Value()
$ws.Range("B15").Value = $ws.Range("B14").Value()
$ws.Range("A14").Value = $ws.Range("A13").Value()
$ws.Range("B14").Value = $ws.Range("B13").Value()
$ws.Range("A13").Value = $ws.Range("A12").Value()
$ws.Range("B13").Value = $ws.Range("B12").Value()

# 2. New "Jurisdiction" row, inserted at row 12.
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

# 3. Field edits.
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-11-22T12:33:30-06:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"
